# Add a "Save" column (H) to the s_vals sheet:
#  - H1 header cell "Save", formatted like the other header cells (bold/border/centered)
#  - H2 data cell with value 0
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new header + value first.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0

# Match the header formatting used by the rest of row 1 (copy G1's style onto H1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
